$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Yr." and "Title" values for several entries were still sitting in columns
# D/E (an older layout) while the rest of the table already used columns F/G for
# those two fields. Consolidate everything onto F (Yr.) / G (Title) so the whole
# list uses the same two columns, then clear out the now-unused D/E cells.
for ($r = 9; $r -le 33; $r++) {
  $dCell = $ws.Cells.Item($r, 4)
  $eCell = $ws.Cells.Item($r, 5)
  $fCell = $ws.Cells.Item($r, 6)
  $gCell = $ws.Cells.Item($r, 7)

  $dVal = $dCell.Value()
  $eVal = $eCell.Value()

  if ($dVal -ne $null) {
    $fCell.Value = $dVal
    $dCell.Clear()
  }
  if ($eVal -ne $null) {
    $gCell.Value = $eVal
    $eCell.Clear()
  }
}

# Now that every row is consistently laid out (call # / author / title-author /
# .. / Yr. / Title), sort the list by call number (column A) so the blank
# spacer rows collapse to the bottom and the catalogued books sit in one
# contiguous, alphabetically-ordered block.
$sortRange = $ws.Range("A9:I33")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A9:A33"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4163
$ws.Sort.Apply()

$ws.Range("F20").Select()
